# edit.ps1 - Apply "Added Universal Design to the DesignDocument" changes
# Word COM-interop (iron_native) script

$d = $word.ActiveDocument

# Useful characters
$cr      = [char]13        # paragraph mark, used inside Find/Replace "Replace" strings
$lq      = [char]0x2018    # left single quote (not used below but handy)
$rq      = [char]0x2019    # right single quote / apostrophe used throughout doc
$ldq     = [char]0x201C    # left double quote
$rdq     = [char]0x201D    # right double quote
$arrow   = [char]0x2192    # RIGHTWARDS ARROW

# ---------------------------------------------------------------------------
# 1) Seat-map sentence: describe seat symbols with quoted "0"/"X" instead of
#    "represented by x's and o's"
# ---------------------------------------------------------------------------
$old1 = "represented by x" + $rq + "s and o" + $rq + "s"
$new1 = "occupied seats are represented by " + $ldq + "0" + $rdq + ", while unoccupied seats are labeled by " + $ldq + "X" + $rdq
$d.Content.Find.Execute($old1, $false, $false, $false, $false, $false, $true, 1, $false, $new1, 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Split the "Users can have exclusive memberships..." paragraph in two,
#    inserting a new sentence about loyalty points for Customers.
# ---------------------------------------------------------------------------
$old2 = "non-member. A user can be an administrator, "
$new2 = "non-member. Customers also have the ability to accumulate loyalty points from purchasing flights. " + $cr + "A user can be an administrator, "
$d.Content.Find.Execute($old2, $false, $false, $false, $false, $false, $true, 1, $false, $new2, 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Two new blank paragraphs right after "Phase 2 refactoring + code smells"
# ---------------------------------------------------------------------------
$heading = $d.Content.Find.Execute("Phase 2 refactoring + code smells", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$phaseHeadingPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13) -eq "Phase 2 refactoring + code smells") {
        $phaseHeadingPara = $i
        break
    }
}

$insertPos = $d.Paragraphs.Item($phaseHeadingPara).Range.End - 1
$ins = $d.Range($insertPos, $insertPos)
$ins.InsertParagraphAfter()
$d.Paragraphs.Item($phaseHeadingPara + 1).Range.Style = $d.Styles.Item("Normal")

$insertPos2 = $d.Paragraphs.Item($phaseHeadingPara + 1).Range.End - 1
$ins2 = $d.Range($insertPos2, $insertPos2)
$ins2.InsertParagraphAfter()
$d.Paragraphs.Item($phaseHeadingPara + 2).Range.Style = $d.Styles.Item("Normal")

# ---------------------------------------------------------------------------
# 4) New "Universal Design" body content, inserted right after the
#    placeholder paragraph and before "Personal Contribution".
# ---------------------------------------------------------------------------
$p3 = "We followed Principle 3(Simple and Intuitive Use) within the design of our menus by grouping similar information together; e.g the User Menu allows Users to view/change any information related to their personal information (viewing/changing birthdates, payment methods, etc.). Furthermore, all the User prompts and messages are basic/simple to understand, even with a limited English background. Throughout the use of the program, the User will be given simple feedback with task completion (i.e. " + $ldq + "Account has been added!" + $rdq + " after registering, " + $ldq + "Your Payment has been saved" + $rdq + " after adding a payment method). There is also feedback if the User enters an invalid input (e.g " + $ldq + "This Seat is already Occupied. Please Try Again" + $rdq + " if the User wants to purchase an occupied seat). Feedback helps with User intuition: the User knows immediately if their inputs are valid, and thus allows better understanding of the program and removes complexity. "

$p4 = "We also followed Principle 4(Perceptible Information) throughout the program.  Within the flightMenu, essential information about the flights (flight number, destination) are uniquely presented from the rest of the console to maximize visibility of the information. By presenting flight information this way, it is easier for the User to choose a flight number and purchase a ticket, instead of the flight number being cluttered with other less important information. "

$p5 = "The program also follows Principle 5(Tolerance for Error). As mentioned before, there is feedback given when a User performs an action. To prevent User from accidently/purposely crashing and disrupting the program with invalid input, we implemented regex expression matchers within the console class. This " + $ldq + "Fail Safe" + $rdq + " feature prevents the program from progressing (and potentially crashing) unless the User enters a valid input that matches the pattern. This can be seen when the User registers, where the User must enter a valid birth month that must be in the format MM, or an email that contains @).  "

$p6 = "One way we could have improved upon our program to accommodate Universal design is to have the program run automatically (perhaps like an app), instead of the User having to navigate to  " + $ldq + "MainApplication" + $rdq + " in the SRC folder to run.  This would better accommodate Users with little experience with java/IntelliJ. "

$old4 = "A description about how your program does or COULD follow the 7 principles of universal design (from the ethics lecture)"
$new4 = $old4 + $cr + $cr + $p3 + $cr + $p4 + $cr + $p5 + $cr + $p6 + $cr
$d.Content.Find.Execute($old4, $false, $false, $false, $false, $false, $true, 1, $false, $new4, 2) | Out-Null

# Find the paragraph index of the placeholder text now (still same text, paragraph unchanged)
$descIdx = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13) -eq $old4) {
        $descIdx = $i
        break
    }
}
# The very next paragraph (descIdx+1) is the new blank paragraph -> reset to Normal style.
$d.Paragraphs.Item($descIdx + 1).Range.Style = $d.Styles.Item("Normal")

# The "We followed Principle 3..." paragraph should carry the lastRenderedPageBreak
# that used to sit on "Personal Contribution"; also remove it from "Personal Contribution".
$pcIdx = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13) -eq "Personal Contribution (for phase 2) ") {
        $pcIdx = $i
        break
    }
}

# Blank paragraph right before "Personal Contribution" heading (after p6)
if ($pcIdx -ne $null) {
    $d.Paragraphs.Item($pcIdx - 1).Range.Style = $d.Styles.Item("Normal")
}

# ---------------------------------------------------------------------------
# 5) Daniel's personal contribution bullet gets expanded.
# ---------------------------------------------------------------------------
$old5 = "Daniel:"
$new5 = "Daniel: GenerateFlights, CancelAllFlights, ClearAllUsers. Repackaged Tests into appropriate test folders and wrote tests for nearly all methods in all UseCases. "
$d.Content.Find.Execute($old5, $false, $false, $false, $false, $false, $true, 1, $false, $new5, 2) | Out-Null
